# Update countries & provincias Spain
#
# 1) "Chequia" moves ahead of "Emiratos Arabes Unidos" in the country list
#    (new daily figures for Chequia; Emiratos Arabes Unidos keeps its old
#    figures but drops one row down).
# 2) "Timor Oriental" moves ahead of "Suazilandia" / "Botsuana"
#    (new daily figures for Timor Oriental; Suazilandia and Botsuana keep
#    their old figures but drop one row down each).
# 3) Daily refresh of a couple of per-country counters (Alemania, Rusia).
# 4) The "last updated" timestamp banner is bumped forward by 30 minutes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Alemania (row 8) : refreshed counters -----------------------------
$ws.Range("B8").Value = 145743
$ws.Range("C8").Value = 1
$ws.Range("E8").Value = 49601

# --- 2) Rusia (row 13) : refreshed counters -------------------------------
$ws.Range("B13").Value = 47121
$ws.Range("C13").Value = 4268
$ws.Range("D13").Value = 3446
$ws.Range("E13").Value = 43270
$ws.Range("G13").Value = 44
$ws.Range("H13").Value = 405

# --- 3) Chequia jumps ahead of Emiratos Arabes Unidos (rows 38-39) -------
# Row 38 now holds Chequia's fresh data.
$ws.Range("A38").Value = "Chequia"
$ws.Range("B38").Value = 6787
$ws.Range("C38").Value = 41
$ws.Range("D38").Value = 1311
$ws.Range("E38").Value = 5288
$ws.Range("F38").Value = 84
$ws.Range("G38").Value = 2
$ws.Range("H38").Value = 188

# Row 39 now holds Emiratos Arabes Unidos, keeping its previous data.
$ws.Range("A39").Value = "Emiratos Arabes Unidos"
$ws.Range("B39").Value = 6781
$ws.Range("C39").Value = 0
$ws.Range("D39").Value = 1286
$ws.Range("E39").Value = 5454
$ws.Range("F39").Value = 1
$ws.Range("G39").Value = 0
$ws.Range("H39").Value = 41

# --- 4) Timor Oriental jumps ahead of Suazilandia / Botsuana -------------
# (rows 178-180)
# Row 178 now holds Timor Oriental's fresh data.
$ws.Range("A178").Value = "Timor Oriental"
$ws.Range("B178").Value = 22
$ws.Range("C178").Value = 3
$ws.Range("D178").Value = 1
$ws.Range("E178").Value = 21
$ws.Range("F178").Value = 0
$ws.Range("G178").Value = 0
$ws.Range("H178").Value = 0

# Row 179 now holds Suazilandia, keeping its previous data.
$ws.Range("A179").Value = "Suazilandia"
$ws.Range("B179").Value = 22
$ws.Range("C179").Value = 0
$ws.Range("D179").Value = 8
$ws.Range("E179").Value = 13
$ws.Range("F179").Value = 0
$ws.Range("G179").Value = 0
$ws.Range("H179").Value = 1

# Row 180 now holds Botsuana, keeping its previous data.
$ws.Range("A180").Value = "Botsuana"
$ws.Range("B180").Value = 20
$ws.Range("C180").Value = 0
$ws.Range("D180").Value = 0
$ws.Range("E180").Value = 19
$ws.Range("F180").Value = 0
$ws.Range("G180").Value = 0
$ws.Range("H180").Value = 1

# Row 181 (Laos) is unaffected - it already holds the correct data.

# --- 5) Bump the "last updated" banner ------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 20 de Abril de 2020 a las 09:52"
